$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.181445333333333
$ws.Range("H2").Value = 6.544335999999999
$ws.Range("I2").Value = 0.1058843243701343
$ws.Range("J2").Value = 0.1058843243701343
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 118.0470123333333
$ws.Range("N2").Value = 354.141037
$ws.Range("O2").Value = 0.4657216250363638
$ws.Range("P2").Value = 0.4657216250363638
$ws.Range("Q2").Value = 257.5131041684924
$ws.Range("R2").Value = 2317.617937516432
$ws.Range("S2").Value = 0.04931261961153638
$ws.Range("T2").Value = 0.04931261961153639

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.181445333333333
$ws.Range("H3").Value = 6.544335999999999
$ws.Range("I3").Value = 0.1058843243701343
$ws.Range("J3").Value = 0.1058843243701343
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 59.36586533333332
$ws.Range("N3").Value = 178.097596
$ws.Range("O3").Value = 0.2342114953037475
$ws.Range("P3").Value = 0.2342114953037476
$ws.Range("Q3").Value = 129.5033898906951
$ws.Range("R3").Value = 1165.530509016256
$ws.Range("S3").Value = 0.02479932593995618
$ws.Range("T3").Value = 0.02479932593995619

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.181445333333333
$ws.Range("H4").Value = 6.544335999999999
$ws.Range("I4").Value = 0.1058843243701343
$ws.Range("J4").Value = 0.1058843243701343
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 56.84506433333333
$ws.Range("N4").Value = 170.535193
$ws.Range("O4").Value = 0.2242663767030476
$ws.Range("P4").Value = 0.2242663767030477
$ws.Range("Q4").Value = 124.0044003129831
$ws.Range("R4").Value = 1116.039602816848
$ws.Range("S4").Value = 0.02374629377614022
$ws.Range("T4").Value = 0.02374629377614022

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.181445333333333
$ws.Range("H5").Value = 6.544335999999999
$ws.Range("I5").Value = 0.1058843243701343
$ws.Range("J5").Value = 0.1058843243701343
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 19.21324333333333
$ws.Range("N5").Value = 57.63973
$ws.Range("O5").Value = 0.07580050295684103
$ws.Range("P5").Value = 0.07580050295684104
$ws.Range("Q5").Value = 41.91264000769777
$ws.Range("R5").Value = 377.21376006928
$ws.Range("S5").Value = 0.008026085042501477
$ws.Range("T5").Value = 0.008026085042501479

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 8.793934999999999
$ws.Range("H6").Value = 26.381805
$ws.Range("I6").Value = 0.426845381730038
$ws.Range("J6").Value = 0.426845381730038
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 118.0470123333333
$ws.Range("N6").Value = 354.141037
$ws.Range("O6").Value = 0.4657216250363638
$ws.Range("P6").Value = 0.4657216250363638
$ws.Range("Q6").Value = 1038.097753403531
$ws.Range("R6").Value = 9342.879780631785
$ws.Range("S6").Value = 0.1987911248185803
$ws.Range("T6").Value = 0.1987911248185803

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 8.793934999999999
$ws.Range("H7").Value = 26.381805
$ws.Range("I7").Value = 0.426845381730038
$ws.Range("J7").Value = 0.426845381730038
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 59.36586533333332
$ws.Range("N7").Value = 178.097596
$ws.Range("O7").Value = 0.2342114953037475
$ws.Range("P7").Value = 0.2342114953037476
$ws.Range("Q7").Value = 522.0595609600865
$ws.Range("R7").Value = 4698.53604864078
$ws.Range("S7").Value = 0.09997209511849112
$ws.Range("T7").Value = 0.09997209511849114

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.793934999999999
$ws.Range("H8").Value = 26.381805
$ws.Range("I8").Value = 0.426845381730038
$ws.Range("J8").Value = 0.426845381730038
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 56.84506433333333
$ws.Range("N8").Value = 170.535193
$ws.Range("O8").Value = 0.2242663767030476
$ws.Range("P8").Value = 0.2242663767030477
$ws.Range("Q8").Value = 499.8918008181516
$ws.Range("R8").Value = 4499.026207363365
$ws.Range("S8").Value = 0.09572706717302487
$ws.Range("T8").Value = 0.09572706717302489

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.793934999999999
$ws.Range("H9").Value = 26.381805
$ws.Range("I9").Value = 0.426845381730038
$ws.Range("J9").Value = 0.426845381730038
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 19.21324333333333
$ws.Range("N9").Value = 57.63973
$ws.Range("O9").Value = 0.07580050295684103
$ws.Range("P9").Value = 0.07580050295684104
$ws.Range("Q9").Value = 168.9600130125167
$ws.Range("R9").Value = 1520.64011711265
$ws.Range("S9").Value = 0.03235509461994168
$ws.Range("T9").Value = 0.03235509461994169

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.198723333333334
$ws.Range("H10").Value = 18.59617
$ws.Range("I10").Value = 0.3008774146563012
$ws.Range("J10").Value = 0.3008774146563012
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 118.0470123333333
$ws.Range("N10").Value = 354.141037
$ws.Range("O10").Value = 0.4657216250363638
$ws.Range("P10").Value = 0.4657216250363638
$ws.Range("Q10").Value = 731.7407697809211
$ws.Range("R10").Value = 6585.66692802829
$ws.Range("S10").Value = 0.1401251184904725
$ws.Range("T10").Value = 0.1401251184904725

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 6.198723333333334
$ws.Range("H11").Value = 18.59617
$ws.Range("I11").Value = 0.3008774146563012
$ws.Range("J11").Value = 0.3008774146563012
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 59.36586533333332
$ws.Range("N11").Value = 178.097596
$ws.Range("O11").Value = 0.2342114953037475
$ws.Range("P11").Value = 0.2342114953037476
$ws.Range("Q11").Value = 367.9925746452577
$ws.Range("R11").Value = 3311.93317180732
$ws.Range("S11").Value = 0.070468949189778
$ws.Range("T11").Value = 0.07046894918977802

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 6.198723333333334
$ws.Range("H12").Value = 18.59617
$ws.Range("I12").Value = 0.3008774146563012
$ws.Range("J12").Value = 0.3008774146563012
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 56.84506433333333
$ws.Range("N12").Value = 170.535193
$ws.Range("O12").Value = 0.2242663767030476
$ws.Range("P12").Value = 0.2242663767030477
$ws.Range("Q12").Value = 352.3668266678678
$ws.Range("R12").Value = 3171.30144001081
$ws.Range("S12").Value = 0.06747668761674912
$ws.Range("T12").Value = 0.06747668761674913

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 6.198723333333334
$ws.Range("H13").Value = 18.59617
$ws.Range("I13").Value = 0.3008774146563012
$ws.Range("J13").Value = 0.3008774146563012
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 19.21324333333333
$ws.Range("N13").Value = 57.63973
$ws.Range("O13").Value = 0.07580050295684103
$ws.Range("P13").Value = 0.07580050295684104
$ws.Range("Q13").Value = 119.0975797593445
$ws.Range("R13").Value = 1071.8782178341
$ws.Range("S13").Value = 0.02280665935930165
$ws.Range("T13").Value = 0.02280665935930165

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.428052000000001
$ws.Range("H14").Value = 10.284156
$ws.Range("I14").Value = 0.1663928792435264
$ws.Range("J14").Value = 0.1663928792435264
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 118.0470123333333
$ws.Range("N14").Value = 354.141037
$ws.Range("O14").Value = 0.4657216250363638
$ws.Range("P14").Value = 0.4657216250363638
$ws.Range("Q14").Value = 404.671296723308
$ws.Range("R14").Value = 3642.041670509772
$ws.Range("S14").Value = 0.07749276211577456
$ws.Range("T14").Value = 0.07749276211577458

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.428052000000001
$ws.Range("H15").Value = 10.284156
$ws.Range("I15").Value = 0.1663928792435264
$ws.Range("J15").Value = 0.1663928792435264
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 59.36586533333332
$ws.Range("N15").Value = 178.097596
$ws.Range("O15").Value = 0.2342114953037475
$ws.Range("P15").Value = 0.2342114953037476
$ws.Range("Q15").Value = 203.509273387664
$ws.Range("R15").Value = 1831.583460488976
$ws.Range("S15").Value = 0.03897112505552222
$ws.Range("T15").Value = 0.03897112505552223

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.428052000000001
$ws.Range("H16").Value = 10.284156
$ws.Range("I16").Value = 0.1663928792435264
$ws.Range("J16").Value = 0.1663928792435264
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 56.84506433333333
$ws.Range("N16").Value = 170.535193
$ws.Range("O16").Value = 0.2242663767030476
$ws.Range("P16").Value = 0.2242663767030477
$ws.Range("Q16").Value = 194.867836478012
$ws.Range("R16").Value = 1753.810528302108
$ws.Range("S16").Value = 0.03731632813713341
$ws.Range("T16").Value = 0.03731632813713342

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.428052000000001
$ws.Range("H17").Value = 10.284156
$ws.Range("I17").Value = 0.1663928792435264
$ws.Range("J17").Value = 0.1663928792435264
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 19.21324333333333
$ws.Range("N17").Value = 57.63973
$ws.Range("O17").Value = 0.07580050295684103
$ws.Range("P17").Value = 0.07580050295684104
$ws.Range("Q17").Value = 65.86399723532001
$ws.Range("R17").Value = 592.7759751178801
$ws.Range("S17").Value = 0.01261266393509622
$ws.Range("T17").Value = 0.01261266393509622
